$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three "checked developer mode" cells (C10, C20, C31) are updated to
# reflect the new "CHROME developer mode" wording used for the added
# exploratory-testing note.
$ws.Range("C10").Value = "/ CHROME developer mode ->  network"
$ws.Range("C20").Value = "/ CHROME developer mode ->  network"
$ws.Range("C31").Value = "/ CHROME developer mode ->  network"

# Window / view state captured by the author when the workbook was re-saved:
# zoomed out from 110% to 80% and moved the selection to C24.
$ws.Range("C24").Select()
$excel.ActiveWindow.Zoom = 80
